$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated roster data (player, positions, team) for rows 2-19
$data = @(
    @("Isaiah Collier", "PG,SG", "Utah Jazz"),
    @("De'Aaron Fox", "PG,SG", "San Antonio Spurs"),
    @("Luka Doncic", "PG,SG", "Los Angeles Lakers"),
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Scottie Barnes", "PG,SG,SF,PF", "Toronto Raptors"),
    @("Zach Collins", "PF,C", "Chicago Bulls"),
    @("Evan Mobley", "PG,SG", "Cleveland Cavaliers"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Keldon Johnson", "SG,SF,PF", "San Antonio Spurs"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Jonathan Kuminga", "SF,PF", "Golden State Warriors"),
    @("P.J. Washington", "SF,PF", "Dallas Mavericks")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
